# Case and Fatality Demographics Data Updated (2021-10-15)
# Updates the three "Fatalities by ..." sheets with refreshed counts.

$wb = $excel.ActiveWorkbook

# --- Fatalities by Age Group -------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")

$wsAge.Range("B2").Value  = 13
$wsAge.Range("B3").Value  = 18
$wsAge.Range("B4").Value  = 68
$wsAge.Range("B5").Value  = 551
$wsAge.Range("B6").Value  = 1808
$wsAge.Range("B7").Value  = 4568
$wsAge.Range("B8").Value  = 8783
$wsAge.Range("B9").Value  = 6755
$wsAge.Range("B10").Value = 8106
$wsAge.Range("B11").Value = 8687
$wsAge.Range("B12").Value = 8303
$wsAge.Range("B13").Value = 19817

# Total row switches from a SUM formula to a hard-coded total value.
$wsAge.Range("B15").Value = 67477

# --- Fatalities by Gender -----------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")

$wsGender.Range("B2").Value = 28183
$wsGender.Range("B3").Value = 39293

# --- Fatalities by Race-Ethnicity ---------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

$wsRace.Range("B2").Value = 1262
$wsRace.Range("B3").Value = 6990
$wsRace.Range("B4").Value = 29670
$wsRace.Range("B5").Value = 394
$wsRace.Range("B6").Value = 29121
$wsRace.Range("B7").Value = 40

# --- Restore the selections / active sheet seen in the saved workbook --------
# Select Race-Ethnicity and Gender sheets first so their selection marks are
# updated without leaving them as the final active sheet.
[void]$wsRace.Activate()
$wsRace.Range("G21").Select() | Out-Null

[void]$wsGender.Activate()
$wsGender.Range("D13").Select() | Out-Null

# Fatalities by Age Group ends up as the active (selected) tab.
[void]$wsAge.Activate()
$wsAge.Range("H10").Select() | Out-Null
